$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row for 2026/01/18 (append after the last existing row, 68).
# Column A must stay plain text "2026/01/18" (not get auto-converted to a
# date serial number), so we enter it with a leading apostrophe to force
# text, then clear the resulting "quote prefix" formatting flag and
# reapply the same centered alignment used by the rest of the table so
# the cell ends up sharing the existing cell style instead of minting a
# new one.
$ws.Range("A69").Value = "'2026/01/18"
$ws.Range("A69").ClearFormats()

$ws.Range("B69").Value = "逃离鸭科夫"
$ws.Range("C69").Value = 1146

$ws.Range("A69:C69").HorizontalAlignment = -4108
$ws.Range("A69:C69").VerticalAlignment = -4108
